# Edit the "Examples & Info" sheet:
#  - Row 5 becomes the "category" row (previously row 6's content), with
#    B5 keeping its existing "sample" value.
#  - Row 6 becomes the "regex" row (previously row 7's content).
#  - The old row 7 is removed entirely (rows shift up), shrinking the used
#    range from A1:BJ7 to A1:BJ6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Examples & Info")

# --- Row 5: label + category values ---------------------------------------
$ws.Range("A5").Value = "category"
# B5 already contains "sample" - leave it as-is (restate for clarity/safety).
$ws.Range("B5:R5").Value = "sample"
$ws.Range("S5:AM5").Value = "experiment"
$ws.Range("AN5:AX5").Value = "sequencing"
$ws.Range("AY5:BJ5").Value = "other"

# --- Row 6: label + regex values -------------------------------------------
$ws.Range("A6").Value = "regex"

# Clear out the old "category" values that used to live in row 6 (B..Y),
# they are blank in the regex row.
$ws.Range("B6:Y6").ClearContents()

$ws.Range("Z6").Value = "[+-]?([0-9]*[.])?[0-9]+"
$ws.Range("AA6").Value = "[+-]?([0-9]*[.])?[0-9]+"
$ws.Range("AB6").Value = "[0-9]*"
$ws.Range("AC6").Value = "[0-9]*"
$ws.Range("AD6").Value = "[+-]?([0-9]*[.])?[0-9]+"
$ws.Range("AE6").Value = "[0-9]*"
$ws.Range("AF6").Value = "[+-]?([0-9]*[.])?[0-9]+"
$ws.Range("AG6:AH6").ClearContents()
$ws.Range("AI6").Value = "[ATGC]*"
$ws.Range("AJ6").ClearContents()
$ws.Range("AK6").Value = "[ATGC, ]*"
$ws.Range("AL6:AN6").ClearContents()
$ws.Range("AO6").Value = "[0-9]*"
$ws.Range("AP6").Value = "[A-Z]{2}[0-9]{4}[0-9]{2}[0-9]{2}"
$ws.Range("AQ6").Value = "[0-9]*"
$ws.Range("AR6").Value = "[0-9]{5,}"
$ws.Range("AS6:AU6").ClearContents()
$ws.Range("AV6").Value = "[0-9]{1,2}"
$ws.Range("AW6").Value = "[A-Z]"
$ws.Range("AX6").Value = "[0-9]{1,2}[A-Z]"

# Clear out the old "other" values that used to live in row 6 (AY..BJ),
# they are blank in the regex row.
$ws.Range("AY6:BJ6").ClearContents()

# --- Remove the now-duplicated old row 7 (regex data already moved up) ----
$ws.Rows("7").Delete()
